# Add a new "Greener scenario" land-use / emission comparison sheet
# (compile_lcArea_Green) after the existing compile_lcArea sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("compile_lcArea")

# Insert the new worksheet right after the original one, then rename it.
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "compile_lcArea_Green"

# --- Data -----------------------------------------------------------------
# Row 1 is the header: ID, Land_use_type, 2017_ha .. 2030_ha
# Rows 2-23 are the 22 land-use categories with yearly areas (ha).
$rows = @(
    @("ID","Land_use_type","2017_ha","2018_ha","2019_ha","2020_ha","2021_ha","2022_ha","2023_ha","2024_ha","2025_ha","2026_ha","2027_ha","2028_ha","2029_ha","2030_ha"),
    @(1,"Hutan Lahan Kering Primer",580931,580932,580948,580954,580961,580971,580981,580993,581000,581007,581008,581009,581008,581011),
    @(2,"Hutan Lahan Kering Sekunder",289816,291556,293204,294841,296452,297813,299130,300349,301415,302366,303277,304112,304909,305577),
    @(3,"Hutan Rawa Primer",132451,132519,132578,132633,132683,132732,132773,132819,132865,132913,132952,132989,133031,133066),
    @(4,"Hutan Rawa Sekunder",46117,47799,48935,49821,50483,51147,51592,52082,52397,52638,52909,53177,53405,53589),
    @(5,"Hutan Mangrove Primer",805,805,805,805,805,805,805,805,805,805,805,805,805,805),
    @(6,"Hutan Mangrove Sekunder",4095,4131,4158,4183,4210,4240,4272,4301,4331,4352,4369,4390,4409,4424),
    @(7,"Hutan Tanaman",240037,227663,220199,215492,212109,208711,205998,203884,202420,201073,199989,198965,198134,197282),
    @(8,"Semak/Belukar",445750,525344,595960,652318,697784,734554,758878,761648,760831,755408,755820,754589,754729,753455),
    @(9,"Belukar Rawa",208972,201669,195233,190723,187554,185074,183548,182129,181197,180317,179820,179095,178850,178001),
    @(10,"Savana/Padang Rumput",92,92,92,92,92,92,92,92,92,92,92,92,92,92),
    @(11,"Perkebunan",971225,1005990,1037350,1065482,1090864,1113844,1134115,1152914,1169050,1177086,1180765,1183927,1185524,1184392),
    @(12,"Pertanian Lahan Kering",98692,94558,92074,90592,89668,89235,89828,90489,90784,90948,91174,91453,91734,92030),
    @(13,"Pertanian Lahan Kering + Semak",1433560,1330873,1236890,1156903,1089073,1029728,983410,959276,941246,935378,927912,922205,916606,915315),
    @(14,"Transmigrasi",21839,21841,21844,21848,21852,21857,21863,21869,21875,21881,21886,21892,21898,21904),
    @(15,"Sawah",16925,16980,17132,17319,17517,17712,17876,18022,18173,18336,18504,18664,18825,18998),
    @(16,"Tambak",200,197,197,197,197,197,197,197,197,197,197,197,197,197),
    @(17,"Tanah Terbuka",233977,237058,237216,235812,233219,232585,231818,230980,229466,228441,226972,226220,225047,224628),
    @(18,"Pertambangan",14568,15569,16482,17363,18234,19045,19820,20560,21275,21977,22666,23348,24015,24664),
    @(19,"Permukiman",91209,95745,100082,104058,107735,111205,114605,118245,122287,126542,130690,134728,138686,142520),
    @(20,"Tubuh Air",43564,43564,43564,43564,43564,43564,43564,43564,43564,43564,43564,43564,43564,43564),
    @(21,"Rawa",15141,15081,15023,14966,14910,14855,14801,14748,14696,14645,14595,14545,14498,14452),
    @(22,"Pelabuhan Udara/Laut",81,81,81,81,81,81,81,81,81,81,81,81,81,81)
)

$nRows = $rows.Count
$nCols = $rows[0].Count
$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}

$startCell = $newSheet.Cells.Item(1,1)
$endCell = $newSheet.Cells.Item($nRows, $nCols)
$dataRange = $newSheet.Range($startCell, $endCell)
$dataRange.Value = $arr

# --- Column widths ----------------------------------------------------------
# Target widths (as stored in OOXML): col A = 3, col B = 27.21875
$newSheet.Columns.Item(1).ColumnWidth = 2.1666666666666665
$newSheet.Columns.Item(2).ColumnWidth = 26.385416666666668

# --- Page setup ---------------------------------------------------------
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# --- Selection on the new sheet (matches the saved view state) -----------
[void]$newSheet.Range("C24:P24").Select()
